# PROS-9738 - CCRU - New POS 2019 KPIs
# Adds 6 new atomic-name-update rows (38-43) to the "Update Atomic Names" sheet,
# each with its KPI Set / KPI Name / Atomic Name Old / Atomic Name New values
# and the usual CONCATENATE(...) UPDATE-statement helper formula in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 38: PoS 2019 - MT Supermarket - REG / SSD Availability
#         Schweppes Tonic - 0.5L -> Schweppes Mojito - 0.33L
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = "PoS 2019 - MT Supermarket - REG"

$ws.Range("B38").Value = "SSD Availability"
$ws.Range("B38").VerticalAlignment = -4108

$ws.Range("C38").Value = "Schweppes Tonic - 0.5L"

$ws.Range("D38").Value = "Schweppes Mojito - 0.33L"
$ws.Range("D38").VerticalAlignment = -4108
$ws.Range("D38").Interior.Color = 5296274

$ws.Range("E38").Formula = '=CONCATENATE("UPDATE `static`.atomic_kpi a JOIN `static`.kpi k ON k.pk=a.kpi_fk JOIN `static`.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name=''",D38,"'', a.description=''",D38,"'', a.display_text=''",D38,"''  WHERE s.name=''",A38,"'' AND k.display_text=''",B38,"'' AND a.name=''",C38,"'';")'

# ---------------------------------------------------------------------------
# Row 39: PoS 2019 - FT NS - CAP / SSD Availability
#         Sprite - 0.25L Slim -> Schweppes Bitter Lemon - 0.33L
# ---------------------------------------------------------------------------
$ws.Range("A39").Value = "PoS 2019 - FT NS - CAP"
$ws.Range("B39").Value = "SSD Availability"
$ws.Range("C39").Value = "Sprite - 0.25L Slim"

$ws.Range("D39").Value = "Schweppes Bitter Lemon - 0.33L"
$ws.Range("D39").VerticalAlignment = -4108
$ws.Range("D39").Interior.Color = 5296274

$ws.Range("E39").Formula = '=CONCATENATE("UPDATE `static`.atomic_kpi a JOIN `static`.kpi k ON k.pk=a.kpi_fk JOIN `static`.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name=''",D39,"'', a.description=''",D39,"'', a.display_text=''",D39,"''  WHERE s.name=''",A39,"'' AND k.display_text=''",B39,"'' AND a.name=''",C39,"'';")'

# ---------------------------------------------------------------------------
# Row 40: PoS 2019 - FT NS - CAP / Energy Availability
#         Burn Apple Kiwi - 0.5L -> Monster Green - 0.5L
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = "PoS 2019 - FT NS - CAP"
$ws.Range("B40").Value = "Energy Availability"
$ws.Range("C40").Value = "Burn Apple Kiwi - 0.5L"

$ws.Range("D40").Value = "Monster Green - 0.5L"
$ws.Range("D40").VerticalAlignment = -4108
$ws.Range("D40").Interior.Color = 5296274

$ws.Range("E40").Formula = '=CONCATENATE("UPDATE `static`.atomic_kpi a JOIN `static`.kpi k ON k.pk=a.kpi_fk JOIN `static`.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name=''",D40,"'', a.description=''",D40,"'', a.display_text=''",D40,"''  WHERE s.name=''",A40,"'' AND k.display_text=''",B40,"'' AND a.name=''",C40,"'';")'

# ---------------------------------------------------------------------------
# Row 41: PoS 2019 - FT NS - REG / SSD Availability
#         Sprite - 0.25L Slim -> Schweppes Bitter Lemon - 0.33L
# ---------------------------------------------------------------------------
$ws.Range("A41").Value = "PoS 2019 - FT NS - REG"
$ws.Range("B41").Value = "SSD Availability"
$ws.Range("C41").Value = "Sprite - 0.25L Slim"

$ws.Range("D41").Value = "Schweppes Bitter Lemon - 0.33L"
$ws.Range("D41").VerticalAlignment = -4108
$ws.Range("D41").Interior.Color = 5296274

$ws.Range("E41").Formula = '=CONCATENATE("UPDATE `static`.atomic_kpi a JOIN `static`.kpi k ON k.pk=a.kpi_fk JOIN `static`.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name=''",D41,"'', a.description=''",D41,"'', a.display_text=''",D41,"''  WHERE s.name=''",A41,"'' AND k.display_text=''",B41,"'' AND a.name=''",C41,"'';")'

# ---------------------------------------------------------------------------
# Row 42: PoS 2019 - IC FastFood / Juice Availability
#         Dobriy - Apple - 0.33L -> Dobriy - Apple - 0.33L/Rich Apple - 0.3L
# ---------------------------------------------------------------------------
$ws.Range("A42").Value = "PoS 2019 - IC FastFood"
$ws.Range("B42").Value = "Juice Availability"
$ws.Range("C42").Value = "Dobriy - Apple - 0.33L"

$ws.Range("D42").Value = "Dobriy - Apple - 0.33L/Rich Apple - 0.3L "
$ws.Range("D42").VerticalAlignment = -4108
$ws.Range("D42").Interior.Color = 5296274
$ws.Range("D42").Borders.LineStyle = -4115

$ws.Range("E42").Formula = '=CONCATENATE("UPDATE `static`.atomic_kpi a JOIN `static`.kpi k ON k.pk=a.kpi_fk JOIN `static`.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name=''",D42,"'', a.description=''",D42,"'', a.display_text=''",D42,"''  WHERE s.name=''",A42,"'' AND k.display_text=''",B42,"'' AND a.name=''",C42,"'';")'

# ---------------------------------------------------------------------------
# Row 43: PoS 2019 - IC FastFood / Juice Availability
#         Dobriy - Multifruit - 0.33L -> Dobriy - Multifruit - 0.33L/Rich Orange 0.3L
# ---------------------------------------------------------------------------
$ws.Range("A43").Value = "PoS 2019 - IC FastFood"
$ws.Range("B43").Value = "Juice Availability"
$ws.Range("C43").Value = "Dobriy - Multifruit - 0.33L"

$ws.Range("D43").Value = "Dobriy - Multifruit - 0.33L/Rich Orange 0.3L"
$ws.Range("D43").VerticalAlignment = -4108
$ws.Range("D43").Interior.Color = 5296274
$ws.Range("D43").Borders.LineStyle = -4115

$ws.Range("E43").Formula = '=CONCATENATE("UPDATE `static`.atomic_kpi a JOIN `static`.kpi k ON k.pk=a.kpi_fk JOIN `static`.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name=''",D43,"'', a.description=''",D43,"'', a.display_text=''",D43,"''  WHERE s.name=''",A43,"'' AND k.display_text=''",B43,"'' AND a.name=''",C43,"'';")'

# ---------------------------------------------------------------------------
# Restore the active selection to A20, matching the saved view state.
# ---------------------------------------------------------------------------
[void]$ws.Range("A20").Select()

Write-Output "KPI atomic-name rows 38-43 added"
